$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell A2 contained "Il était une fois" - update it to add the trailing period.
$ws.Range("A2").Value = "Il était une fois."
